$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'260.65"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'1.69%"
$ws.Range("E2").Style = "Normal"
$ws.Range("G2").Value = "'14"
$ws.Range("G2").Style = "Normal"
$ws.Range("D3").Value = "'27.22"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'2.61%"
$ws.Range("E3").Style = "Normal"
$ws.Range("G3").Value = "'14"
$ws.Range("G3").Style = "Normal"
$ws.Range("D4").Value = "'4.696"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'0.50%"
$ws.Range("E4").Style = "Normal"
$ws.Range("G4").Value = "'14"
$ws.Range("G4").Style = "Normal"
$ws.Range("D5").Value = "'0.06125"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'3.41%"
$ws.Range("E5").Style = "Normal"
$ws.Range("G5").Value = "'14"
$ws.Range("G5").Style = "Normal"
$ws.Range("E6").Value = "'0.71%"
$ws.Range("E6").Style = "Normal"
$ws.Range("G6").Value = "'14"
$ws.Range("G6").Style = "Normal"
$ws.Range("D7").Value = "'0.8512"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'-0.36%"
$ws.Range("E7").Style = "Normal"
$ws.Range("G7").Value = "'14"
$ws.Range("G7").Style = "Normal"
$ws.Range("D8").Value = "'0.9219"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'0.88%"
$ws.Range("E8").Style = "Normal"
$ws.Range("G8").Value = "'14"
$ws.Range("G8").Style = "Normal"
$ws.Range("E9").Value = "'2.17%"
$ws.Range("E9").Style = "Normal"
$ws.Range("G9").Value = "'14"
$ws.Range("G9").Style = "Normal"
$ws.Range("D10").Value = "'0.04633"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'5.58%"
$ws.Range("E10").Style = "Normal"
$ws.Range("G10").Value = "'14"
$ws.Range("G10").Style = "Normal"
$ws.Range("D11").Value = "'0.07093"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'1.35%"
$ws.Range("E11").Style = "Normal"
$ws.Range("G11").Value = "'14"
$ws.Range("G11").Style = "Normal"
$ws.Range("D12").Value = "'0.03085"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'1.76%"
$ws.Range("E12").Style = "Normal"
$ws.Range("G12").Value = "'14"
$ws.Range("G12").Style = "Normal"
$ws.Range("D13").Value = "'0.09061"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'-0.46%"
$ws.Range("E13").Style = "Normal"
$ws.Range("G13").Value = "'14"
$ws.Range("G13").Style = "Normal"
$ws.Range("D14").Value = "'0.001533"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'0.47%"
$ws.Range("E14").Style = "Normal"
$ws.Range("G14").Value = "'14"
$ws.Range("G14").Style = "Normal"
$ws.Range("D15").Value = "'0.0006086"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-94.10%"
$ws.Range("E15").Style = "Normal"
$ws.Range("G15").Value = "'14"
$ws.Range("G15").Style = "Normal"
$ws.Range("D16").Value = "'0.006106"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'0.33%"
$ws.Range("E16").Style = "Normal"
$ws.Range("G16").Value = "'14"
$ws.Range("G16").Style = "Normal"
$ws.Range("D17").Value = "'3.450"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.72%"
$ws.Range("E17").Style = "Normal"
$ws.Range("G17").Value = "'14"
$ws.Range("G17").Style = "Normal"
$ws.Range("D18").Value = "'3.154"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.52%"
$ws.Range("E18").Style = "Normal"
$ws.Range("G18").Value = "'14"
$ws.Range("G18").Style = "Normal"
$ws.Range("D19").Value = "'2.163"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'0.60%"
$ws.Range("E19").Style = "Normal"
$ws.Range("G19").Value = "'14"
$ws.Range("G19").Style = "Normal"
$ws.Range("E20").Value = "'0.78%"
$ws.Range("E20").Style = "Normal"
$ws.Range("G20").Value = "'14"
$ws.Range("G20").Style = "Normal"
$ws.Range("D21").Value = "'0.1304"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'1.54%"
$ws.Range("E21").Style = "Normal"
$ws.Range("G21").Value = "'14"
$ws.Range("G21").Style = "Normal"
$ws.Range("D22").Value = "'4.096"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'6.12%"
$ws.Range("E22").Style = "Normal"
$ws.Range("G22").Value = "'14"
$ws.Range("G22").Style = "Normal"
$ws.Range("D23").Value = "'0.04228"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.84%"
$ws.Range("E23").Style = "Normal"
$ws.Range("G23").Value = "'14"
$ws.Range("G23").Style = "Normal"
$ws.Range("D24").Value = "'0.001216"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'0.03%"
$ws.Range("E24").Style = "Normal"
$ws.Range("G24").Value = "'14"
$ws.Range("G24").Style = "Normal"
$ws.Range("D25").Value = "'0.003789"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-18.56%"
$ws.Range("E25").Style = "Normal"
$ws.Range("G25").Value = "'14"
$ws.Range("G25").Style = "Normal"
$ws.Range("E26").Value = "'0.16%"
$ws.Range("E26").Style = "Normal"
$ws.Range("G26").Value = "'14"
$ws.Range("G26").Style = "Normal"
$ws.Range("E27").Value = "'-8.03%"
$ws.Range("E27").Style = "Normal"
$ws.Range("G27").Value = "'14"
$ws.Range("G27").Style = "Normal"
$ws.Range("G28").Value = "'14"
$ws.Range("G28").Style = "Normal"
$ws.Range("G29").Value = "'14"
$ws.Range("G29").Style = "Normal"
$ws.Range("G30").Value = "'14"
$ws.Range("G30").Style = "Normal"
$ws.Range("G31").Value = "'14"
$ws.Range("G31").Style = "Normal"
$ws.Range("G32").Value = "'14"
$ws.Range("G32").Style = "Normal"
$ws.Range("G33").Value = "'14"
$ws.Range("G33").Style = "Normal"
$ws.Range("G34").Value = "'14"
$ws.Range("G34").Style = "Normal"
$ws.Range("G35").Value = "'14"
$ws.Range("G35").Style = "Normal"
$ws.Range("G36").Value = "'14"
$ws.Range("G36").Style = "Normal"
$ws.Range("G37").Value = "'14"
$ws.Range("G37").Style = "Normal"
$ws.Range("G38").Value = "'14"
$ws.Range("G38").Style = "Normal"
$ws.Range("G39").Value = "'14"
$ws.Range("G39").Style = "Normal"
$ws.Range("D40").Value = "'0.03870"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'2.13%"
$ws.Range("E40").Style = "Normal"
$ws.Range("G40").Value = "'14"
$ws.Range("G40").Style = "Normal"
$ws.Range("D41").Value = "'0.1110"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'0.88%"
$ws.Range("E41").Style = "Normal"
$ws.Range("G41").Value = "'14"
$ws.Range("G41").Style = "Normal"
$ws.Range("D42").Value = "'0.004077"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-34.99%"
$ws.Range("E42").Style = "Normal"
$ws.Range("G42").Value = "'14"
$ws.Range("G42").Style = "Normal"
$ws.Range("D43").Value = "'0.01633"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'12.64%"
$ws.Range("E43").Style = "Normal"
$ws.Range("G43").Value = "'14"
$ws.Range("G43").Style = "Normal"
$ws.Range("E44").Value = "'-3.94%"
$ws.Range("E44").Style = "Normal"
$ws.Range("G44").Value = "'14"
$ws.Range("G44").Style = "Normal"
$ws.Range("D45").Value = "'0.00005160"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'0.76%"
$ws.Range("E45").Style = "Normal"
$ws.Range("G45").Value = "'14"
$ws.Range("G45").Style = "Normal"
$ws.Range("E46").Value = "'0.16%"
$ws.Range("E46").Style = "Normal"
$ws.Range("G46").Value = "'14"
$ws.Range("G46").Style = "Normal"
$ws.Range("E47").Value = "'8.15%"
$ws.Range("E47").Style = "Normal"
$ws.Range("G47").Value = "'14"
$ws.Range("G47").Style = "Normal"
$ws.Range("D48").Value = "'0.1355"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'-43.74%"
$ws.Range("E48").Style = "Normal"
$ws.Range("G48").Value = "'14"
$ws.Range("G48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002101"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'0.16%"
$ws.Range("E49").Style = "Normal"
$ws.Range("G49").Value = "'14"
$ws.Range("G49").Style = "Normal"
$ws.Range("E50").Value = "'0.16%"
$ws.Range("E50").Style = "Normal"
$ws.Range("G50").Value = "'14"
$ws.Range("G50").Style = "Normal"
$ws.Range("G51").Value = "'14"
$ws.Range("G51").Style = "Normal"
